$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MRNA")

# Update "Change in payables and accrued liability" row (row 7) values
$ws.Range("B7").Value = 69000000.0
$ws.Range("C7").Value = 58000000.0
$ws.Range("D7").Value = 56543000.0
$ws.Range("E7").Value = 42304000.0
$ws.Range("F7").Value = -1253000.0
